# Applies the "article rewrite" edit to bmark_combined.xlsx:
#  - Processing!B2:D3 formulas now divide by 1,000,000 (µs/ns -> seconds-ish scale)
#  - Processing!B5:B10 lose their (self-referential / now redundant) formulas
#  - Processing sheet formatting is rebuilt (drops a stray border flag that was
#    baked into several cellXfs) and columns B:D are widened
#  - Selection/active-tab moves from Results to Processing

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Results")
$ws2 = $wb.Worksheets.Item("Processing")

# --- Processing: rebuild number formatting from a clean slate -------------
# (ClearFormats drops the old, unused "applyBorder" cruft that every cellXf
# in this sheet was carrying, so the re-applied formats below start fresh.)
$ws2.Range("A1:D16").ClearFormats()

# Plain 2-decimal cells (avg / deviation / accel ratios that aren't %).
$ws2.Range("B2:D3").NumberFormat  = "0.00"
$ws2.Range("B5:D9").NumberFormat  = "0.00"
$ws2.Range("B11:D15").NumberFormat = "0.00"

# Row 4 ("rel. error") keeps the Percent cell style.
$ws2.Range("B4:D4").Style = "Percent"
$ws2.Range("B4:D4").NumberFormat = "0.00%"

# Rows 10 and 16 ("max./mid. rel. accel" ratios) are plain percent cells.
$ws2.Range("B10:D10").NumberFormat = "0.00%"
$ws2.Range("B16:D16").NumberFormat = "0.00%"

# --- Processing: rescale the raw-time formulas to millions -----------------
$ws2.Range("B2").Formula = "=AVERAGE(Results!B:B)/1000000"
$ws2.Range("C2").Formula = "=AVERAGE(Results!C:C)/1000000"
$ws2.Range("D2").Formula = "=AVERAGE(Results!D:D)/1000000"

$ws2.Range("B3").Formula = "=_xlfn.STDEV.S(Results!B:B)/1000000"
$ws2.Range("C3").Formula = "=_xlfn.STDEV.S(Results!C:C)/1000000"
$ws2.Range("D3").Formula = "=_xlfn.STDEV.S(Results!D:D)/1000000"

# --- Processing: column B's self-referential accel/error formulas (rows
# 5-10) are no longer needed now that the table only has two real
# benchmarks to compare against the baseline; drop them but keep the cells
# (and their number format) in place.
$ws2.Range("B5:B10").ClearContents()

# --- Processing: widen the data columns ------------------------------------
$ws2.Range("B1:D16").ColumnWidth = 9.6

# --- Selection / active sheet ----------------------------------------------
[void]$ws1.Range("F17").Select()
[void]$ws2.Activate()
[void]$ws2.Range("D3").Select()

# --- Window geometry (best effort) -----------------------------------------
$excel.ActiveWindow.Left   = 4788
$excel.ActiveWindow.Top    = 1008
$excel.ActiveWindow.Width  = 17280
$excel.ActiveWindow.Height = 8964
